# "ngoc anh update tien do"
# Fill in the newly-added "Xem Hinh Cham Cong" row (row 36) of the
# "Cham cong va bao bieu" section with owner / progress / note info.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B36").Value = "Xem Hình Chấm Công"
$ws.Range("D36").Value = "Ngọc Anh"
$ws.Range("E36").Value = "xong front-end - Xong backend"
$ws.Range("F36").Value = "Cần chú ý vì hình được lưu thẳng vào db dạng binary"

# Move the selection to where the author left off editing.
$ws.Range("E49").Select()
